$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 5
$ws.Range("E2").Value = 1
$ws.Range("A3").Value = 8
$ws.Range("E3").Value = 3
$ws.Range("A4").Value = 11
$ws.Range("E4").Value = 5
$ws.Range("A5").Value = 15
$ws.Range("E5").Value = 8
$ws.Range("A6").Value = 16
$ws.Range("E6").Value = 9
$ws.Range("A7").Value = 20
$ws.Range("E7").Value = 12
$ws.Range("A8").Value = 23
$ws.Range("E8").Value = 14
$ws.Range("A9").Value = 24
$ws.Range("E9").Value = 15
$ws.Range("A10").Value = 26
$ws.Range("E10").Value = 17
$ws.Range("A11").Value = 30
$ws.Range("E11").Value = 19
$ws.Range("A12").Value = 33
$ws.Range("E12").Value = 21
$ws.Range("A13").Value = 36
$ws.Range("E13").Value = 23
$ws.Range("A14").Value = 39
$ws.Range("E14").Value = 25
$ws.Range("A15").Value = 41
$ws.Range("E15").Value = 26
$ws.Range("A16").Value = 14
$ws.Range("E16").Value = 10
$ws.Range("A17").Value = 5
$ws.Range("E17").Value = 6
$ws.Range("A18").Value = 21
$ws.Range("E18").Value = 16
$ws.Range("A19").Value = 18
$ws.Range("E19").Value = 13
$ws.Range("A20").Value = 31
$ws.Range("E20").Value = 22
$ws.Range("A21").Value = 32
$ws.Range("E21").Value = 24
$ws.Range("A22").Value = 27
$ws.Range("E22").Value = 20
$ws.Range("A23").Value = 10
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = '2023-06-21'
$ws.Range("B23").NumberFormat = "General"
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").Value = '19:00'
$ws.Range("E23").Value = 11
$ws.Range("F23").Value = 'Wed'
$ws.Range("H23").Value = 'L'
$ws.Range("I23").Value = 0
$ws.Range("L23").Value = 1.1
$ws.Range("M23").Value = 1.8
$ws.Range("N23").Value = 67
$ws.Range("O23").Value = 28494
$ws.Range("P23").Value = 12
$ws.Range("Q23").Value = 3
$ws.Range("R23").Value = 25
$ws.Range("S23").Value = 0
$ws.Range("T23").Value = 0
$ws.Range("X23").Value = 1.1
$ws.Range("Y23").Value = 0.09
$ws.Range("Z23").Value = -1.1
$ws.Range("AA23").Value = -1.1
$ws.Range("AF23").Value = 1
$ws.Range("AG23").Value = 0
$ws.Range("AH23").Value = 9581
$ws.Range("AI23").Value = 3140
$ws.Range("AK23").Value = 1
$ws.Range("AM23").Value = 10
$ws.Range("AN23").Value = 53
$ws.Range("AO23").Value = 12
$ws.Range("AP23").Value = 5
$ws.Range("AQ23").Value = 51
$ws.Range("AR23").Value = 0
$ws.Range("AS23").Value = 7
$ws.Range("AT23").Value = 27
$ws.Range("AU23").Value = 17
$ws.Range("AV23").Value = 23
$ws.Range("AW23").Value = 0
$ws.Range("AX23").Value = 9
$ws.Range("AY23").Value = 1
$ws.Range("AZ23").Value = 10
$ws.Range("BA23").Value = 3
$ws.Range("BB23").Value = 5
$ws.Range("BD23").Value = 'Cruzeiro'
$ws.Range("A24").Value = 6
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = '2023-05-08'
$ws.Range("B24").NumberFormat = "General"
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value = '20:00'
$ws.Range("E24").Value = 4
$ws.Range("F24").Value = 'Mon'
$ws.Range("H24").Value = 'D'
$ws.Range("I24").Value = 1
$ws.Range("L24").Value = 1.5
$ws.Range("M24").Value = 0.8
$ws.Range("N24").Value = 49
$ws.Range("O24").Value = 36512
$ws.Range("P24").Value = 21
$ws.Range("Q24").Value = 8
$ws.Range("R24").Value = 38.1
$ws.Range("S24").Value = 0.05
$ws.Range("T24").Value = 0.13
$ws.Range("X24").Value = 1.5
$ws.Range("Y24").Value = 0.07
$ws.Range("Z24").Value = -0.5
$ws.Range("AA24").Value = -0.5
$ws.Range("AF24").Value = 0.8
$ws.Range("AG24").Value = -0.2
$ws.Range("AH24").Value = 8153
$ws.Range("AI24").Value = 2923
$ws.Range("AK24").Value = 1.3
$ws.Range("AM24").Value = 15
$ws.Range("AN24").Value = 40
$ws.Range("AO24").Value = 5
$ws.Range("AP24").Value = 2
$ws.Range("AQ24").Value = 37
$ws.Range("AR24").Value = 3
$ws.Range("AS24").Value = 3
$ws.Range("AT24").Value = 17
$ws.Range("AU24").Value = 11
$ws.Range("AV24").Value = 35
$ws.Range("AW24").Value = 2
$ws.Range("AX24").Value = 12
$ws.Range("AY24").Value = 9
$ws.Range("AZ24").Value = 4
$ws.Range("BA24").Value = 0
$ws.Range("BB24").Value = 8
$ws.Range("BD24").Value = 'Corinthians'
$ws.Range("A25").Value = 24
$ws.Range("E25").Value = 18
$ws.Range("A26").Value = 1
$ws.Range("E26").Value = 2
$ws.Range("A27").Value = 9
$ws.Range("E27").Value = 7
